$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 3150
$ws.Range("I20").Value = 1300
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 1300
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -1070
$ws.Range("N20").Value = -5460
$ws.Range("H28").Value = 310.85715
$ws.Range("I28").Value = 126
$ws.Range("J28").Value = 449.5
$ws.Range("K28").Value = 126
$ws.Range("L28").Value = 449.5
$ws.Range("M28").Value = 359
$ws.Range("N28").Value = -1419.5
$ws.Range("H35").Value = 3150
$ws.Range("I35").Value = 1300
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 1300
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -921
$ws.Range("N35").Value = -5758
$ws.Range("H41").Value = 228.3
$ws.Range("I41").Value = 193
$ws.Range("J41").Value = 263.6
$ws.Range("K41").Value = 193
$ws.Range("L41").Value = 263.6
$ws.Range("M41").Value = 247
$ws.Range("N41").Value = -1143.6
$ws.Range("H64").Value = 5249.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 5249.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 5249.5
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -5745.5
$ws.Range("H67").Value = 5249.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 5249.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 5249.5
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -6965.5
$ws.Range("H107").Value = 575.6
$ws.Range("I107").Value = 479.2
$ws.Range("J107").Value = 672
$ws.Range("K107").Value = 479.2
$ws.Range("L107").Value = 672
$ws.Range("M107").Value = 1440.8
$ws.Range("N107").Value = -4512
$ws.Range("H113").Value = 7494.926
$ws.Range("I113").Value = 5999.2
$ws.Range("J113").Value = 8374.764999999999
$ws.Range("K113").Value = 5999.2
$ws.Range("L113").Value = 8374.764999999999
$ws.Range("M113").Value = -2745.2
$ws.Range("N113").Value = -14882.765
$ws.Range("H115").Value = 124.666664
$ws.Range("I115").Value = 124.666664
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 373.999992
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 1193.000008
$ws.Range("H137").Value = 2251
$ws.Range("I137").Value = 2002
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 6006
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -3456
$ws.Range("N137").Value = -12600

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6000
$ws.Range("I61").Value = 6000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5788
$ws.Range("H88").Value = 2927.6
$ws.Range("I88").Value = 249
$ws.Range("J88").Value = 3225.2222
$ws.Range("K88").Value = 249
$ws.Range("L88").Value = 3225.2222
$ws.Range("M88").Value = 157
$ws.Range("N88").Value = -4037.2222
$ws.Range("H91").Value = 2927.6
$ws.Range("I91").Value = 249
$ws.Range("J91").Value = 3225.2222
$ws.Range("K91").Value = 249
$ws.Range("L91").Value = 3225.2222
$ws.Range("M91").Value = 1155
$ws.Range("N91").Value = -6033.2222
$ws.Range("H97").Value = 772.5454999999999
$ws.Range("I97").Value = 772.5454999999999
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 772.5454999999999
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -276.5454999999999
$ws.Range("H136").Value = 6000
$ws.Range("I136").Value = 6000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 18000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15450

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -875
$ws.Range("H86").Value = 6166.6665
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 4250
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 4250
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = -6496
$ws.Range("H89").Value = 6166.6665
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 4250
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 21250
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = -32482
$ws.Range("H94").Value = 1269.1428
$ws.Range("I94").Value = 1214
$ws.Range("J94").Value = 1600
$ws.Range("K94").Value = 1214
$ws.Range("L94").Value = 1600
$ws.Range("M94").Value = -763
$ws.Range("N94").Value = -2502

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3250
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3250
$ws.Range("N62").Value = -4498
$ws.Range("H65").Value = 3250
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 16250
$ws.Range("N65").Value = -22490
$ws.Range("H88").Value = 32999.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 32999.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 32999.5
$ws.Range("N88").Value = -33811.5
$ws.Range("H91").Value = 32999.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 32999.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 32999.5
$ws.Range("N91").Value = -35807.5

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 287500.4
$ws.Range("I11").Value = 220000
$ws.Range("J11").Value = 355000.8
$ws.Range("K11").Value = 220000
$ws.Range("L11").Value = 355000.8
$ws.Range("M11").Value = -219861
$ws.Range("N11").Value = -355278.8
$ws.Range("H19").Value = 22499.5
$ws.Range("I19").Value = 24999
$ws.Range("J19").Value = 20000
$ws.Range("K19").Value = 24999
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = -24711
$ws.Range("N19").Value = -20576
$ws.Range("H80").Value = 3342.2856
$ws.Range("I80").Value = 3179.2
$ws.Range("J80").Value = 3750
$ws.Range("K80").Value = 3179.2
$ws.Range("L80").Value = 3750
$ws.Range("M80").Value = -2181.2
$ws.Range("N80").Value = -5746
$ws.Range("H83").Value = 3342.2856
$ws.Range("I83").Value = 3179.2
$ws.Range("J83").Value = 3750
$ws.Range("K83").Value = 15896
$ws.Range("L83").Value = 18750
$ws.Range("M83").Value = -10904
$ws.Range("N83").Value = -28734
$ws.Range("H123").Value = 35001
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 35001
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 35001
$ws.Range("N123").Value = -39901

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18716.572
$ws.Range("I7").Value = 19104.75
$ws.Range("J7").Value = 18199
$ws.Range("K7").Value = 19104.75
$ws.Range("L7").Value = 18199
$ws.Range("M7").Value = -18992.75
$ws.Range("N7").Value = -18423
$ws.Range("H22").Value = 652.13336
$ws.Range("I22").Value = 482.5
$ws.Range("J22").Value = 1330.6666
$ws.Range("K22").Value = 482.5
$ws.Range("L22").Value = 1330.6666
$ws.Range("M22").Value = -187.5
$ws.Range("N22").Value = -1920.6666
$ws.Range("H27").Value = 652.13336
$ws.Range("I27").Value = 482.5
$ws.Range("J27").Value = 1330.6666
$ws.Range("K27").Value = 482.5
$ws.Range("L27").Value = 1330.6666
$ws.Range("M27").Value = -375.5
$ws.Range("N27").Value = -1544.6666
$ws.Range("H126").Value = 18716.572
$ws.Range("I126").Value = 19104.75
$ws.Range("J126").Value = 18199
$ws.Range("K126").Value = 57314.25
$ws.Range("L126").Value = 54597
$ws.Range("M126").Value = -54844.25
$ws.Range("N126").Value = -59537

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2330.3809
$ws.Range("I107").Value = 2759.0833
$ws.Range("J107").Value = 1758.7778
$ws.Range("K107").Value = 8277.249899999999
$ws.Range("L107").Value = 5276.3334
$ws.Range("M107").Value = -6357.249899999999
$ws.Range("N107").Value = -9116.3334

Write-Host "Applied all Marilith_Profits cell updates."